$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 49.52770149253723
$ws.Range("B2").Value = 33.71200000000002
$ws.Range("C2").Value = 65.59999999999998
$ws.Range("A3").Value = 2.228716417910446
$ws.Range("B3").Value = 0.8680000000000001
$ws.Range("C3").Value = 3.831999999999998
$ws.Range("A4").Value = 3.93725373134327
$ws.Range("B4").Value = 1.867999999999999
$ws.Range("C4").Value = 6.480000000000001
$ws.Range("A5").Value = 2.818169154228847
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 4.743999999999998
$ws.Range("A6").Value = 12.81954228855721
$ws.Range("B6").Value = 6.940000000000003
$ws.Range("C6").Value = 19.05199999999999
$ws.Range("A7").Value = 42.92593034825866
$ws.Range("B7").Value = 26.28400000000001
$ws.Range("C7").Value = 59.09200000000003
$ws.Range("A8").Value = 71.20316417910445
$ws.Range("B8").Value = 54.96400000000001
$ws.Range("C8").Value = 84.50799999999995
$ws.Range("A9").Value = 11.15862686567164
$ws.Range("B9").Value = 5.591999999999996
$ws.Range("C9").Value = 17.52800000000001
$ws.Range("A10").Value = 72.41432835820885
$ws.Range("B10").Value = 56.69600000000002
$ws.Range("C10").Value = 85.33999999999996
$ws.Range("A11").Value = 64.06805970149244
$ws.Range("B11").Value = 46.41199999999999
$ws.Range("C11").Value = 78.23600000000002
$ws.Range("A12").Value = 32.87078606965166
$ws.Range("B12").Value = 19.136
$ws.Range("C12").Value = 46.068
$ws.Range("A13").Value = 73.72121393034816
$ws.Range("B13").Value = 56.83200000000002
$ws.Range("C13").Value = 85.84399999999998
$ws.Range("A14").Value = 69.96911442786059
$ws.Range("B14").Value = 53.672
$ws.Range("C14").Value = 82.20799999999997
$ws.Range("A15").Value = 3.409711442786061
$ws.Range("B15").Value = 1.588
$ws.Range("C15").Value = 5.504
$ws.Range("A16").Value = 12.90794029850745
$ws.Range("B16").Value = 7.084000000000004
$ws.Range("C16").Value = 19.17999999999999
$ws.Range("A17").Value = 32.21554228855717
$ws.Range("B17").Value = 19.58800000000001
$ws.Range("C17").Value = 44.99599999999997
$ws.Range("A18").Value = 3.93725373134327
$ws.Range("B18").Value = 1.867999999999999
$ws.Range("C18").Value = 6.480000000000001
$ws.Range("A19").Value = 66.32865671641787
$ws.Range("B19").Value = 48.89200000000002
$ws.Range("C19").Value = 79.42800000000001
$ws.Range("A20").Value = 14.88171144278607
$ws.Range("B20").Value = 8.815999999999997
$ws.Range("C20").Value = 21.67600000000001
$ws.Range("A21").Value = 5.877333333333326
$ws.Range("B21").Value = 3.04
$ws.Range("C21").Value = 9.612000000000005
$ws.Range("A22").Value = 74.04875621890537
$ws.Range("B22").Value = 57.96800000000004
$ws.Range("C22").Value = 86.77200000000005
$ws.Range("A23").Value = 3.641870646766156
$ws.Range("B23").Value = 1.568
$ws.Range("C23").Value = 6.183999999999998
$ws.Range("A24").Value = 12.93397014925372
$ws.Range("B24").Value = 7.064000000000004
$ws.Range("C24").Value = 19.31999999999998
$ws.Range("A25").Value = 9.970885572139288
$ws.Range("B25").Value = 5.459999999999999
$ws.Range("C25").Value = 15.58399999999999
$ws.Range("A26").Value = 5.951641791044763
$ws.Range("B26").Value = 3.044
$ws.Range("C26").Value = 9.408000000000007
$ws.Range("A27").Value = 74.73158208955222
$ws.Range("B27").Value = 58.96400000000004
$ws.Range("C27").Value = 86.676
$ws.Range("A28").Value = 46.79265671641782
$ws.Range("B28").Value = 30.51999999999999
$ws.Range("C28").Value = 61.62
$ws.Range("A29").Value = 3.726805970149245
$ws.Range("B29").Value = 1.671999999999999
$ws.Range("C29").Value = 6.403999999999999
$ws.Range("A30").Value = 30.51824875621877
$ws.Range("B30").Value = 18.47200000000001
$ws.Range("C30").Value = 43.58400000000002
$ws.Range("A31").Value = 69.41054726368156
$ws.Range("B31").Value = 52.02799999999996
$ws.Range("C31").Value = 82.49999999999993
$ws.Range("A32").Value = 70.955462686567
$ws.Range("B32").Value = 53.59200000000003
$ws.Range("C32").Value = 84.15999999999994
$ws.Range("A33").Value = 70.11797014925369
$ws.Range("B33").Value = 54.25999999999996
$ws.Range("C33").Value = 82.94399999999996
$ws.Range("A34").Value = 54.19797014925361
$ws.Range("B34").Value = 37.7
$ws.Range("C34").Value = 71.98
$ws.Range("A35").Value = 13.53418905472635
$ws.Range("B35").Value = 7.559999999999997
$ws.Range("C35").Value = 20.91999999999999
$ws.Range("A36").Value = 43.21424875621883
$ws.Range("B36").Value = 27.364
$ws.Range("C36").Value = 59.48400000000004
$ws.Range("A37").Value = 6.615980099502476
$ws.Range("B37").Value = 3.464000000000002
$ws.Range("C37").Value = 10.308
$ws.Range("A38").Value = 71.42348258706454
$ws.Range("B38").Value = 56.15600000000001
$ws.Range("C38").Value = 84.25599999999999
$ws.Range("A39").Value = 33.63709452736317
$ws.Range("B39").Value = 20.71200000000002
$ws.Range("C39").Value = 48.63199999999996
$ws.Range("A40").Value = 5.022288557213923
$ws.Range("B40").Value = 2.396000000000002
$ws.Range("C40").Value = 7.632000000000004
$ws.Range("A41").Value = 39.59781094527357
$ws.Range("B41").Value = 24.42799999999999
$ws.Range("C41").Value = 55.27999999999999
$ws.Range("A42").Value = 65.75540298507458
$ws.Range("B42").Value = 48.41999999999998
$ws.Range("C42").Value = 79.33599999999997
$ws.Range("A43").Value = 30.51824875621877
$ws.Range("B43").Value = 18.47200000000001
$ws.Range("C43").Value = 43.58400000000002
$ws.Range("A44").Value = 56.45836815920387
$ws.Range("B44").Value = 39.06000000000002
$ws.Range("C44").Value = 71.77199999999998
$ws.Range("A45").Value = 67.29490547263669
$ws.Range("B45").Value = 50.20400000000004
$ws.Range("C45").Value = 81.62000000000006
$ws.Range("A46").Value = 66.99090547263683
$ws.Range("B46").Value = 49.488
$ws.Range("C46").Value = 81.22799999999992
$ws.Range("A47").Value = 74.73962189054706
$ws.Range("B47").Value = 58.09600000000001
$ws.Range("C47").Value = 86.87600000000006
$ws.Range("A48").Value = 32.61691542288551
$ws.Range("B48").Value = 19.36800000000001
$ws.Range("C48").Value = 45.768
$ws.Range("A49").Value = 4.453572139303466
$ws.Range("B49").Value = 2.04
$ws.Range("C49").Value = 7.06
$ws.Range("A50").Value = 6.535621890547255
$ws.Range("B50").Value = 3.416000000000001
$ws.Range("C50").Value = 10.24
$ws.Range("A51").Value = 2.217273631840794
$ws.Range("B51").Value = 0.852
$ws.Range("C51").Value = 3.811999999999998
$ws.Range("A52").Value = 65.19410945273624
$ws.Range("B52").Value = 46.66399999999998
$ws.Range("C52").Value = 79.35599999999994
$ws.Range("A53").Value = 33.71004975124373
$ws.Range("B53").Value = 21.06000000000001
$ws.Range("C53").Value = 48.82399999999996
$ws.Range("A54").Value = 3.38855721393034
$ws.Range("B54").Value = 1.58
$ws.Range("C54").Value = 5.504
$ws.Range("A55").Value = 6.004875621890537
$ws.Range("B55").Value = 3.148
$ws.Range("C55").Value = 9.120000000000006
$ws.Range("A56").Value = 69.53665671641784
$ws.Range("B56").Value = 53.03999999999996
$ws.Range("C56").Value = 82.06399999999994
$ws.Range("A57").Value = 70.93540298507448
$ws.Range("B57").Value = 53.52000000000002
$ws.Range("C57").Value = 84.20799999999994
$ws.Range("A58").Value = 6.502965174129343
$ws.Range("B58").Value = 3.428000000000003
$ws.Range("C58").Value = 10.152
$ws.Range("A59").Value = 61.97032835820882
$ws.Range("B59").Value = 44.94800000000001
$ws.Range("C59").Value = 76.97999999999996
$ws.Range("A60").Value = 3.72507462686566
$ws.Range("B60").Value = 1.54
$ws.Range("C60").Value = 6.380000000000003
$ws.Range("A61").Value = 71.47160199004971
$ws.Range("B61").Value = 55.45199999999997
$ws.Range("C61").Value = 83.59199999999997
$ws.Range("A62").Value = 64.05548258706457
$ws.Range("B62").Value = 46.27199999999999
$ws.Range("C62").Value = 78.30400000000002
$ws.Range("A63").Value = 72.34991044776115
$ws.Range("B63").Value = 56.66799999999999
$ws.Range("C63").Value = 84.88799999999996
$ws.Range("A64").Value = 3.730149253731335
$ws.Range("B64").Value = 1.647999999999999
$ws.Range("C64").Value = 6.319999999999999
$ws.Range("A65").Value = 69.13751243781088
$ws.Range("B65").Value = 52.06400000000002
$ws.Range("C65").Value = 82.65200000000003
$ws.Range("A66").Value = 53.74975124378098
$ws.Range("B66").Value = 36.44400000000002
$ws.Range("C66").Value = 68.32399999999996
$ws.Range("A67").Value = 51.38628855721389
$ws.Range("B67").Value = 34.45199999999999
$ws.Range("C67").Value = 67.82400000000005
$ws.Range("A68").Value = 5.384318407960186
$ws.Range("B68").Value = 2.735999999999999
$ws.Range("C68").Value = 8.588000000000001
$ws.Range("A69").Value = 74.76465671641787
$ws.Range("B69").Value = 58.60799999999998
$ws.Range("C69").Value = 87.21999999999998
$ws.Range("A70").Value = 4.953870646766164
$ws.Range("B70").Value = 2.264000000000001
$ws.Range("C70").Value = 7.499999999999997
$ws.Range("A71").Value = 20.46165174129339
$ws.Range("B71").Value = 10.652
$ws.Range("C71").Value = 31.66799999999999
$ws.Range("A72").Value = 25.58720398009934
$ws.Range("B72").Value = 13.992
$ws.Range("C72").Value = 38.58399999999997
